# priorizadas_semana.xlsx :: df_com_prioridade
# Refresh the prioritisation export: drop the "Proprietário" column and
# reload the ticket rows with the newest crawl (ranks/ages/timestamps).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("df_com_prioridade")
$ws.Activate()

# --- 1. Remove column C ("Proprietário"); D:H shift left to C:G ---------
$ws.Columns.Item(3).Delete()

# --- 2. Header row -------------------------------------------------------
$ws.Range("A1").Value = "ticket"
$ws.Range("B1").Value = "Título"
$ws.Range("C1").Value = "Priorizada"
$ws.Range("D1").Value = "Atualizado em"
$ws.Range("E1").Value = "Data_priorizado"
$ws.Range("F1").Value = "hora_priorizado"
$ws.Range("G1").Value = "Município"

# --- 3. Data rows ----------------------------------------------------------
$ws.Range("A2").Value = 108566
$ws.Range("B2").Value = "[Petrolina-PE][Regulação] - Permitir selecionar lista de Unidades na agenda de serviço"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "1 dia 6 horas atrás"
$ws.Range("E2").Value = "22/03/2024"
$ws.Range("F2").Value = "18:28:05"
$ws.Range("G2").Value = "Petrolina-PE"

$ws.Range("A3").Value = 108902
$ws.Range("B3").Value = "[Petrolina-PE][Ambulatório] Sistema não agenda na data selecionada pela Solicitação de Procedimento de Serviço"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "1 dia 6 horas atrás"
$ws.Range("E3").Value = "22/03/2024"
$ws.Range("F3").Value = "18:28:05"
$ws.Range("G3").Value = "Petrolina-PE"

$ws.Range("A4").Value = 108553
$ws.Range("B4").Value = "[Petrolina-PE][Regulação] Erro na negação de Unidade executante"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "3 horas 5 minutos atrás"
$ws.Range("E4").Value = "22/03/2024"
$ws.Range("F4").Value = "18:28:05"
$ws.Range("G4").Value = "Petrolina-PE"

$ws.Range("A5").Value = 108649
$ws.Range("B5").Value = "[Governador Valadares-MG][Hospital] - Pacientes internados no mesmo leito."
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "8 horas 9 minutos atrás"
$ws.Range("E5").Value = "22/03/2024"
$ws.Range("F5").Value = "18:28:05"
$ws.Range("G5").Value = "Governador Valadares-MG"

$ws.Range("A6").Value = 105678
$ws.Range("B6").Value = "[Volta Redonda-RJ][Hospital] Evolução do atendimento não salva"
$ws.Range("C6").Value = 5
# "11/03/2024" is ambiguous as a date (day<=12) — force text so Excel
# doesn't silently convert it to a date serial number.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "11/03/2024"
$ws.Range("E6").Value = "22/03/2024"
$ws.Range("F6").Value = "18:28:05"
$ws.Range("G6").Value = "Volta Redonda-RJ"

$ws.Range("A7").Value = 105890
$ws.Range("B7").Value = "[Volta Redonda - RJ][Regulação] A função “Negação de procedimentos entre município solicitantes e municípios e unidades executantes” não está bloqueando as unidades parametrizadas"
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = "3 horas 1 minuto atrás"
$ws.Range("E7").Value = "22/03/2024"
$ws.Range("F7").Value = "18:28:05"
$ws.Range("G7").Value = "Volta Redonda - RJ"

$ws.Range("A8").Value = 108729
$ws.Range("B8").Value = "[Volta Redonda-RJ][Ambulatório][Exportação de Produção RAAS] Dados do Tipo de Gestão está incorreto no arquivo RAAS (Gestão Municipal)"
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = "1 dia 6 horas atrás"
$ws.Range("E8").Value = "22/03/2024"
$ws.Range("F8").Value = "18:28:05"
$ws.Range("G8").Value = "Volta Redonda-RJ"

$ws.Range("A9").Value = 108918
$ws.Range("B9").Value = "[Governador Valadares-MG][Backup] - Oficio para disponibilização do backup dos dados da Saúde"
$ws.Range("C9").Value = 7
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "11/03/2024"
$ws.Range("E9").Value = "22/03/2024"
$ws.Range("F9").Value = "18:28:05"
$ws.Range("G9").Value = "Governador Valadares-MG"

$ws.Range("A10").Value = 108332
$ws.Range("B10").Value = "[Volta Redonda-RJ][Almoxarifado][Movimentação Consolidada de Produtos] Valores do relatório divergem"
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = "8 horas 48 minutos atrás"
$ws.Range("E10").Value = "22/03/2024"
$ws.Range("F10").Value = "18:28:05"
$ws.Range("G10").Value = "Volta Redonda-RJ"

$ws.Range("A11").Value = 108529
$ws.Range("B11").Value = "[São Luis- MA] [APP CIDADÃO] Melhoria no APP - REGULAÇÃO"
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = "15/03/2024"
$ws.Range("E11").Value = "22/03/2024"
$ws.Range("F11").Value = "18:28:05"
$ws.Range("G11").Value = "São Luis- MA"
